# Ajuste consumo actualizar proyectos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data (project queue entry refreshed with new project info)
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "0725-TEST"

# C2/F2 hold date-looking text that must stay plain text (not auto-converted
# to a date serial by Excel) - prefix with a single quote to force text entry,
# then restore the Normal style so no stray number-format/quote-prefix mark
# is left behind on the cell.
$ws.Range("C2").Value = "'2025-08-29"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "F:\TEST\0725-TEST"
$ws.Range("F2").Value = "'2025-08-28"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "Sin Dias Rastreos"
$ws.Range("H2").Value = "JEFFERSON BETANCOURT"
